# MVP V2 Scrape Inventory and Apply Purchase Amount at checkout
#
# The "0004" sheet holds an Item/Quantity table (A1:B6). This sorts the
# item list (A2:B6) alphabetically by Item, leaving the header row (row 1)
# untouched, then leaves the selection where the user ended up (G15) —
# mirroring an interactive Data > Sort A→Z pass over the inventory grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data body (A2:B6) ascending by column A, using the worksheet's
# persistent Sort object so the workbook remembers the last-used sort
# (this is what shows up as <sortState> in the saved sheet XML).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A6"))
$ws.Sort.SetRange($ws.Range("A1:B6"))
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply() | Out-Null

# Leave the selection where the user clicked next.
$ws.Range("G15").Select() | Out-Null
